$wb = $excel.ActiveWorkbook
$wsPolitical = $wb.Worksheets.Item("Political")

# Replace the Political sheet's contents (previously a 0,0 (x,y) style
# coordinate matrix) with the first three columns (name, surname, city)
# of the "name" sheet's first 8 rows.
$wsPolitical.Range("A1:D4").ClearContents()

$data = @(
  @("Shraddha", "Kapoor", "Mumbai"),
  @("Aarti", "Devi", "Gaya"),
  @("Puja", "Kumari", "Kodarma"),
  @("Archana", "Tyagi", "delhi"),
  @("Kunal", "Relan", "Pune"),
  @("Mukesh", "Ambani", "London"),
  @("Aparna", "Garimalla", "Houston"),
  @("Steven", "Hensley", "Texas")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsPolitical.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Widen column G on the Political sheet (left over from the user's
# editing session) and move the selection there.
$wsPolitical.Columns.Item(7).ColumnWidth = 32.1

# Make "Political" the active sheet/tab and select the cell the user
# ended up on.
$wsPolitical.Activate()
$wsPolitical.Range("G16").Select()
